# The deck's single slide master points at ppt/theme/theme2.xml, which
# currently carries the "Integral" / "Red Violet" palette, while
# ppt/theme/theme1.xml (only used by the notes master) carries the
# default "Office Theme" palette. The commit swaps the two themes'
# contents, so the presentation-facing theme (theme2.xml) ends up with
# the stock Office colours (name/fontScheme/fmtScheme were already
# identical between the two themes).
#
# Re-colour the live theme's 12-slot colour scheme to the Office
# palette via ThemeColorScheme, which writes straight through to the
# clrScheme of the theme part actually referenced by the slide master.

$p = $ppt.ActivePresentation

$tcs = $p.Slides.Item(1).ThemeColorScheme

# Order matches a:clrScheme child order: dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink.
$tcs.Item(1).RGB  = 0        # dk1      000000
$tcs.Item(2).RGB  = 16777215 # lt1      FFFFFF
$tcs.Item(3).RGB  = 6968388  # dk2      44546A
$tcs.Item(4).RGB  = 15132391 # lt2      E7E6E6
$tcs.Item(5).RGB  = 13998939 # accent1  5B9BD5
$tcs.Item(6).RGB  = 3243501  # accent2  ED7D31
$tcs.Item(7).RGB  = 10855845 # accent3  A5A5A5
$tcs.Item(8).RGB  = 49407    # accent4  FFC000
$tcs.Item(9).RGB  = 12874308 # accent5  4472C4
$tcs.Item(10).RGB = 4697456  # accent6  70AD47
$tcs.Item(11).RGB = 12673797 # hlink    0563C1
$tcs.Item(12).RGB = 7491477  # folHlink 954F72
